$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting (bold font, thin border, centered alignment) from H1
# into the two new header cells I1 and J1, then set their text.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New numeric data for columns I (I0) and J (IF), rows 2-15.
$data = @{
    2  = @(9, 9)
    3  = @(8, 8)
    4  = @(5, 6)
    5  = @(7, 7)
    6  = @(8, 8)
    7  = @(8, 8)
    8  = @(9, 9)
    9  = @(6, 7)
    10 = @(3, 3)
    11 = @(8, 8)
    12 = @(7, 7)
    13 = @(9, 9)
    14 = @(8, 8)
    15 = @(6, 6)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("I$row").Value = $vals[0]
    $ws.Range("J$row").Value = $vals[1]
}
